# Update "想去人数" (F column) counts on the gh-pages data refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 4737
$ws1.Range("F3").Value  = 1885
$ws1.Range("F6").Value  = 3184
$ws1.Range("F8").Value  = 603
$ws1.Range("F9").Value  = 290
$ws1.Range("F10").Value = 655
$ws1.Range("F11").Value = 556
$ws1.Range("F12").Value = 558
$ws1.Range("F15").Value = 1800
$ws1.Range("F16").Value = 1391
$ws1.Range("F18").Value = 1652
$ws1.Range("F20").Value = 132
$ws1.Range("F21").Value = 621
$ws1.Range("F26").Value = 59
$ws1.Range("F27").Value = 114
$ws1.Range("F28").Value = 8
$ws1.Range("F32").Value = 4034
$ws1.Range("F33").Value = 14
$ws1.Range("F36").Value = 1584
$ws1.Range("F38").Value = 1910

# ---- Sheet "演出" ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 29

# ---- Sheet "全部类型" ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 4737
$ws4.Range("F3").Value  = 1885
$ws4.Range("F6").Value  = 3184
$ws4.Range("F8").Value  = 603
$ws4.Range("F9").Value  = 290
$ws4.Range("F10").Value = 656
$ws4.Range("F11").Value = 556
$ws4.Range("F12").Value = 558
$ws4.Range("F13").Value = 29
$ws4.Range("F16").Value = 1800
$ws4.Range("F17").Value = 1391
$ws4.Range("F19").Value = 1652
$ws4.Range("F21").Value = 132
$ws4.Range("F22").Value = 621
$ws4.Range("F27").Value = 59
$ws4.Range("F28").Value = 114
$ws4.Range("F29").Value = 8
$ws4.Range("F33").Value = 4034
$ws4.Range("F35").Value = 14
$ws4.Range("F39").Value = 1584
$ws4.Range("F41").Value = 1910
